$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.844441333333333
$ws.Range("H2").Value = 5.533324
$ws.Range("I2").Value = 0.09360395274144985
$ws.Range("J2").Value = 0.09360395274144986
$ws.Range("M2").Value = 136.544502
$ws.Range("N2").Value = 409.633506
$ws.Range("O2").Value = 0.9681180443787725
$ws.Range("P2").Value = 0.9681180443787725
$ws.Range("Q2").Value = 251.848323328216
$ws.Range("R2").Value = 2266.634909953944
$ws.Range("S2").Value = 0.09061967567417546
$ws.Range("T2").Value = 0.09061967567417548
$ws.Range("G3").Value = 1.844441333333333
$ws.Range("H3").Value = 5.533324
$ws.Range("I3").Value = 0.09360395274144985
$ws.Range("J3").Value = 0.09360395274144986
$ws.Range("O3").Value = 0.002758738216274633
$ws.Range("P3").Value = 0.002758738216274633
$ws.Range("Q3").Value = 0.7176641302208889
$ws.Range("R3").Value = 6.458977171988
$ws.Range("S3").Value = 0.0002582288016222024
$ws.Range("T3").Value = 0.0002582288016222024
$ws.Range("G4").Value = 1.844441333333333
$ws.Range("H4").Value = 5.533324
$ws.Range("I4").Value = 0.09360395274144985
$ws.Range("J4").Value = 0.09360395274144986
$ws.Range("M4").Value = 4.025396333333334
$ws.Range("N4").Value = 12.076189
$ws.Range("O4").Value = 0.0285405766544606
$ws.Range("P4").Value = 0.02854057665446059
$ws.Range("Q4").Value = 7.424607380248446
$ws.Range("R4").Value = 66.82146642223601
$ws.Range("S4").Value = 0.002671510788377857
$ws.Range("T4").Value = 0.002671510788377857
$ws.Range("G5").Value = 1.844441333333333
$ws.Range("H5").Value = 5.533324
$ws.Range("I5").Value = 0.09360395274144985
$ws.Range("J5").Value = 0.09360395274144986
$ws.Range("M5").Value = 0.08217633333333334
$ws.Range("N5").Value = 0.246529
$ws.Range("O5").Value = 0.0005826407504923545
$ws.Range("P5").Value = 0.0005826407504923544
$ws.Range("Q5").Value = 0.1515694258217778
$ws.Range("R5").Value = 1.364124832396
$ws.Range("S5").Value = 0.00005453747727432923
$ws.Range("T5").Value = 0.00005453747727432922
$ws.Range("I6").Value = 0.3174745301946251
$ws.Range("J6").Value = 0.3174745301946252
$ws.Range("M6").Value = 136.544502
$ws.Range("N6").Value = 409.633506
$ws.Range("O6").Value = 0.9681180443787725
$ws.Range("P6").Value = 0.9681180443787725
$ws.Range("Q6").Value = 854.188586990338
$ws.Range("R6").Value = 7687.697282913043
$ws.Range("S6").Value = 0.30735282131209
$ws.Range("T6").Value = 0.3073528213120901
$ws.Range("I7").Value = 0.3174745301946251
$ws.Range("J7").Value = 0.3174745301946252
$ws.Range("O7").Value = 0.002758738216274633
$ws.Range("P7").Value = 0.002758738216274633
$ws.Range("S7").Value = 0.0008758291191417472
$ws.Range("T7").Value = 0.0008758291191417473
$ws.Range("I8").Value = 0.3174745301946251
$ws.Range("J8").Value = 0.3174745301946252
$ws.Range("M8").Value = 4.025396333333334
$ws.Range("N8").Value = 12.076189
$ws.Range("O8").Value = 0.0285405766544606
$ws.Range("P8").Value = 0.02854057665446059
$ws.Range("Q8").Value = 25.18188250484145
$ws.Range("R8").Value = 226.636942543573
$ws.Range("S8").Value = 0.009060906164858565
$ws.Range("T8").Value = 0.009060906164858565
$ws.Range("I9").Value = 0.3174745301946251
$ws.Range("J9").Value = 0.3174745301946252
$ws.Range("M9").Value = 0.08217633333333334
$ws.Range("N9").Value = 0.246529
$ws.Range("O9").Value = 0.0005826407504923545
$ws.Range("P9").Value = 0.0005826407504923544
$ws.Range("Q9").Value = 0.5140747889947779
$ws.Range("R9").Value = 4.626673100953
$ws.Range("S9").Value = 0.0001849735985348041
$ws.Range("T9").Value = 0.0001849735985348041
$ws.Range("G10").Value = 5.721023666666667
$ws.Range("H10").Value = 17.163071
$ws.Range("I10").Value = 0.2903374692647943
$ws.Range("J10").Value = 0.2903374692647943
$ws.Range("M10").Value = 136.544502
$ws.Range("N10").Value = 409.633506
$ws.Range("O10").Value = 0.9681180443787725
$ws.Range("P10").Value = 0.9681180443787725
$ws.Range("Q10").Value = 781.1743274952139
$ws.Range("R10").Value = 7030.568947456925
$ws.Range("S10").Value = 0.2810809429545146
$ws.Range("T10").Value = 0.2810809429545146
$ws.Range("G11").Value = 5.721023666666667
$ws.Range("H11").Value = 17.163071
$ws.Range("I11").Value = 0.2903374692647943
$ws.Range("J11").Value = 0.2903374692647943
$ws.Range("O11").Value = 0.002758738216274633
$ws.Range("P11").Value = 0.002758738216274633
$ws.Range("Q11").Value = 2.226025517597444
$ws.Range("R11").Value = 20.034229658377
$ws.Range("S11").Value = 0.0008009650720772496
$ws.Range("T11").Value = 0.0008009650720772496
$ws.Range("G12").Value = 5.721023666666667
$ws.Range("H12").Value = 17.163071
$ws.Range("I12").Value = 0.2903374692647943
$ws.Range("J12").Value = 0.2903374692647943
$ws.Range("M12").Value = 4.025396333333334
$ws.Range("N12").Value = 12.076189
$ws.Range("O12").Value = 0.0285405766544606
$ws.Range("P12").Value = 0.02854057665446059
$ws.Range("Q12").Value = 23.02938769071323
$ws.Range("R12").Value = 207.264489216419
$ws.Range("S12").Value = 0.008286398797213958
$ws.Range("T12").Value = 0.008286398797213957
$ws.Range("G13").Value = 5.721023666666667
$ws.Range("H13").Value = 17.163071
$ws.Range("I13").Value = 0.2903374692647943
$ws.Range("J13").Value = 0.2903374692647943
$ws.Range("M13").Value = 0.08217633333333334
$ws.Range("N13").Value = 0.246529
$ws.Range("O13").Value = 0.0005826407504923545
$ws.Range("P13").Value = 0.0005826407504923544
$ws.Range("Q13").Value = 0.4701327478398889
$ws.Range("R13").Value = 4.231194730558999
$ws.Range("S13").Value = 0.0001691624409884906
$ws.Range("T13").Value = 0.0001691624409884906
$ws.Range("G14").Value = 5.883520333333333
$ws.Range("H14").Value = 17.650561
$ws.Range("I14").Value = 0.2985840477991308
$ws.Range("J14").Value = 0.2985840477991308
$ws.Range("M14").Value = 136.544502
$ws.Range("N14").Value = 409.633506
$ws.Range("O14").Value = 0.9681180443787725
$ws.Range("P14").Value = 0.9681180443787725
$ws.Range("Q14").Value = 803.3623539218739
$ws.Range("R14").Value = 7230.261185296866
$ws.Range("S14").Value = 0.2890646044379925
$ws.Range("T14").Value = 0.2890646044379925
$ws.Range("G15").Value = 5.883520333333333
$ws.Range("H15").Value = 17.650561
$ws.Range("I15").Value = 0.2985840477991308
$ws.Range("J15").Value = 0.2985840477991308
$ws.Range("O15").Value = 0.002758738216274633
$ws.Range("P15").Value = 0.002758738216274633
$ws.Range("Q15").Value = 2.289252266445222
$ws.Range("R15").Value = 20.603270398007
$ws.Range("S15").Value = 0.0008237152234334339
$ws.Range("T15").Value = 0.0008237152234334339
$ws.Range("G16").Value = 5.883520333333333
$ws.Range("H16").Value = 17.650561
$ws.Range("I16").Value = 0.2985840477991308
$ws.Range("J16").Value = 0.2985840477991308
$ws.Range("M16").Value = 4.025396333333334
$ws.Range("N16").Value = 12.076189
$ws.Range("O16").Value = 0.0285405766544606
$ws.Range("P16").Value = 0.02854057665446059
$ws.Range("Q16").Value = 23.68350117689211
$ws.Range("R16").Value = 213.151510592029
$ws.Range("S16").Value = 0.008521760904010221
$ws.Range("T16").Value = 0.008521760904010221
$ws.Range("G17").Value = 5.883520333333333
$ws.Range("H17").Value = 17.650561
$ws.Range("I17").Value = 0.2985840477991308
$ws.Range("J17").Value = 0.2985840477991308
$ws.Range("M17").Value = 0.08217633333333334
$ws.Range("N17").Value = 0.246529
$ws.Range("O17").Value = 0.0005826407504923545
$ws.Range("P17").Value = 0.0005826407504923544
$ws.Range("Q17").Value = 0.4834861280854444
$ws.Range("R17").Value = 4.351375152769
$ws.Range("S17").Value = 0.0001739672336947306
$ws.Range("T17").Value = 0.0001739672336947306
